# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
#
# The source produced new scrapes for a handful of already-present
# fixtures (so their betting-odds columns B..AC get refreshed / reordered)
# and the four still-unplayed fixtures at the bottom of the sheet were
# dropped (rows 166-169), shrinking the used range from A1:AC169 to
# A1:AC165.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) rows 137-139: cyclic rotation of the B:AC payload -----------------
# new137 <- old139, new138 <- old137, new139 <- old138
# (column A, the running index, stays put on each row)
$row137 = $ws.Range("B137:AC137").Value2
$row138 = $ws.Range("B138:AC138").Value2
$row139 = $ws.Range("B139:AC139").Value2

$ws.Range("B137:AC137").Value2 = $row139
$ws.Range("B138:AC138").Value2 = $row137
$ws.Range("B139:AC139").Value2 = $row138

# --- 2) rows 143-144: simple swap of the B:AC payload ----------------------
$row143 = $ws.Range("B143:AC143").Value2
$row144 = $ws.Range("B144:AC144").Value2

$ws.Range("B143:AC143").Value2 = $row144
$ws.Range("B144:AC144").Value2 = $row143

# --- 3) rows 147-148: simple swap of the B:AC payload ----------------------
$row147 = $ws.Range("B147:AC147").Value2
$row148 = $ws.Range("B148:AC148").Value2

$ws.Range("B147:AC147").Value2 = $row148
$ws.Range("B148:AC148").Value2 = $row147

# --- 4) rows 153-154: simple swap of the B:AC payload ----------------------
$row153 = $ws.Range("B153:AC153").Value2
$row154 = $ws.Range("B154:AC154").Value2

$ws.Range("B153:AC153").Value2 = $row154
$ws.Range("B154:AC154").Value2 = $row153

# --- 5) drop the four not-yet-played fixtures at the bottom ---------------
$ws.Range("A166:AC169").EntireRow.Delete()
